# Auto-generated edit script applying the Phantom_Profits diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3237.375
$ws.Range("I40").Value = 3285.5715
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 3285.5715
$ws.Range("L40").Value = 2900
$ws.Range("M40").Value = -3110.5715
$ws.Range("N40").Value = -3250
$ws.Range("H58").Value = 3772.5625
$ws.Range("J58").Value = 5882.1
$ws.Range("L58").Value = 17646.3
$ws.Range("N58").Value = -17946.3
$ws.Range("H64").Value = 4362.5454
$ws.Range("I64").Value = 9999
$ws.Range("J64").Value = 3798.9
$ws.Range("K64").Value = 9999
$ws.Range("L64").Value = 3798.9
$ws.Range("M64").Value = -9751
$ws.Range("N64").Value = -4294.9
$ws.Range("H67").Value = 4362.5454
$ws.Range("I67").Value = 9999
$ws.Range("J67").Value = 3798.9
$ws.Range("K67").Value = 9999
$ws.Range("L67").Value = 3798.9
$ws.Range("M67").Value = -9141
$ws.Range("N67").Value = -5514.9
$ws.Range("H70").Value = 3173.2666
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 4199.875
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 12599.625
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -13139.625
$ws.Range("H73").Value = 3173.2666
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 4199.875
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 12599.625
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -14471.625
$ws.Range("H76").Value = 5713.7144
$ws.Range("I76").Value = 5713.7144
$ws.Range("K76").Value = 5713.7144
$ws.Range("M76").Value = -5398.7144
$ws.Range("H79").Value = 5713.7144
$ws.Range("I79").Value = 5713.7144
$ws.Range("K79").Value = 5713.7144
$ws.Range("M79").Value = -4621.7144
$ws.Range("H107").Value = 705.75
$ws.Range("I107").Value = 41
$ws.Range("K107").Value = 41
$ws.Range("M107").Value = 1879
$ws.Range("H135").Value = 1491
$ws.Range("I135").Value = 1491
$ws.Range("K135").Value = 13419
$ws.Range("M135").Value = -10884

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21218.217
$ws.Range("I86").Value = 24439.334
$ws.Range("K86").Value = 24439.334
$ws.Range("M86").Value = -23316.334
$ws.Range("H89").Value = 21218.217
$ws.Range("I89").Value = 24439.334
$ws.Range("K89").Value = 122196.67
$ws.Range("M89").Value = -116580.67
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1587
$ws.Range("I16").Value = 1587
$ws.Range("K16").Value = 1587
$ws.Range("M16").Value = -1300
$ws.Range("H113").Value = 1587
$ws.Range("I113").Value = 1587
$ws.Range("K113").Value = 1587
$ws.Range("M113").Value = 583
$ws.Range("H122").Value = 2898
$ws.Range("I122").Value = 2898
$ws.Range("K122").Value = 8694
$ws.Range("M122").Value = -6244
$ws.Range("H132").Value = 18194964
$ws.Range("I132").Value = 22237290
$ws.Range("K132").Value = 66711870
$ws.Range("M132").Value = -66709340

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1516878.6
$ws.Range("I4").Value = 17767
$ws.Range("K4").Value = 53301
$ws.Range("M4").Value = -53189
$ws.Range("H7").Value = 59349.1
$ws.Range("I7").Value = 84324.42999999999
$ws.Range("J7").Value = 1073.3334
$ws.Range("K7").Value = 252973.29
$ws.Range("L7").Value = 3220.0002
$ws.Range("M7").Value = -252861.29
$ws.Range("N7").Value = -3444.0002
$ws.Range("H25").Value = 425
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 350
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1050
$ws.Range("M25").Value = -1331
$ws.Range("N25").Value = -1388
$ws.Range("H30").Value = 425
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 350
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 1050
$ws.Range("M30").Value = -1398
$ws.Range("N30").Value = -1254
$ws.Range("H34").Value = 43909.56
$ws.Range("J34").Value = 47682.176
$ws.Range("L34").Value = 143046.528
$ws.Range("N34").Value = -143214.528
$ws.Range("H115").Value = 794.25
$ws.Range("I115").Value = 809
$ws.Range("K115").Value = 2427
$ws.Range("M115").Value = -1252
$ws.Range("H131").Value = 2202.25
$ws.Range("I131").Value = 1059.8334
$ws.Range("J131").Value = 2691.8572
$ws.Range("K131").Value = 3179.5002
$ws.Range("L131").Value = 8075.571599999999
$ws.Range("M131").Value = 1860.4998
$ws.Range("N131").Value = -18155.5716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 20000
$ws.Range("J4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("N4").Value = -20224

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 100004000
$ws.Range("I22").Value = 4927.6665
$ws.Range("J22").Value = 250002600
$ws.Range("K22").Value = 4927.6665
$ws.Range("L22").Value = 250002600
$ws.Range("M22").Value = -4632.6665
$ws.Range("N22").Value = -250003190
$ws.Range("H27").Value = 100004000
$ws.Range("I27").Value = 4927.6665
$ws.Range("J27").Value = 250002600
$ws.Range("K27").Value = 4927.6665
$ws.Range("L27").Value = 250002600
$ws.Range("M27").Value = -4820.6665
$ws.Range("N27").Value = -250002814
$ws.Range("H68").Value = 1458.875
$ws.Range("I68").Value = 1524.4286
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1524.4286
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -775.4286
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1458.875
$ws.Range("I71").Value = 1524.4286
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 7622.143
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -3878.143
$ws.Range("N71").Value = -12488

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H51").Value = 39969.25
$ws.Range("J51").Value = 30038.5
$ws.Range("L51").Value = 30038.5
$ws.Range("N51").Value = -31058.5
$ws.Range("H52").Value = 20000
$ws.Range("I52").Value = 20000
$ws.Range("K52").Value = 20000
$ws.Range("M52").Value = -19774
$ws.Range("H54").Value = 56049.6
$ws.Range("J54").Value = 56049.6
$ws.Range("L54").Value = 56049.6
$ws.Range("N54").Value = -57089.6
$ws.Range("H55").Value = 4243.2856
$ws.Range("I55").Value = 2930
$ws.Range("J55").Value = 7526.5
$ws.Range("K55").Value = 2930
$ws.Range("L55").Value = 7526.5
$ws.Range("M55").Value = -2653
$ws.Range("N55").Value = -8080.5
$ws.Range("H107").Value = 871.2941
$ws.Range("I107").Value = 850.0526
$ws.Range("J107").Value = 898.2
$ws.Range("K107").Value = 2550.1578
$ws.Range("L107").Value = 2694.6
$ws.Range("M107").Value = -630.1578
$ws.Range("N107").Value = -6534.6
$ws.Range("H122").Value = 2346.3333
$ws.Range("I122").Value = 2346.3333
$ws.Range("K122").Value = 7038.999899999999
$ws.Range("M122").Value = -4588.999899999999
$ws.Range("H136").Value = 9520.154
$ws.Range("I136").Value = 9520.154
$ws.Range("K136").Value = 28560.462
$ws.Range("M136").Value = -26010.462
